$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
# (xlShiftToRight = -4161)
$ws.Columns("D").Insert(-4161)

# Copy number formats / styles from the (now-shifted) old column D, which is
# now column E, into the freshly inserted blank column D so the new column
# matches the look of its neighbours (date format on header rows, number
# format on data rows).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the FY2018 (period ending 2018-12-31) figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 371000
$ws.Range("D9").Value = 250800
$ws.Range("D10").Value = 120200
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 4400
$ws.Range("D17").Value = 323800
$ws.Range("D18").Value = 47200
$ws.Range("D20").Value = -11500
$ws.Range("D21").Value = 58800
$ws.Range("D22").Value = 21200
$ws.Range("D23").Value = 14500
$ws.Range("D24").Value = 6000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 8500
$ws.Range("D27").Value = 9000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 11500
$ws.Range("D33").Value = 9000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 9000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 33000
$ws.Range("D42").Value = 1200
$ws.Range("D43").Value = 148900
$ws.Range("D44").Value = 91800
$ws.Range("D45").Value = 18500
$ws.Range("D46").Value = 293400
$ws.Range("D47").Value = 9800
$ws.Range("D48").Value = 149200
$ws.Range("D49").Value = 32600
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 4900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 489800
$ws.Range("D57").Value = 65500
$ws.Range("D58").Value = 21600
$ws.Range("D59").Value = 35500
$ws.Range("D60").Value = 122600
$ws.Range("D61").Value = 229800
$ws.Range("D62").Value = 4100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 357400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 11800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 132400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 9000
$ws.Range("D83").Value = 23200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -5000
$ws.Range("D91").Value = -13100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -18700
$ws.Range("D96").Value = -2700
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 17000
$ws.Range("D101").Value = -1200
$ws.Range("D102").Value = -7900
